# Add two new instrument model values to the cv_experiment controlled-vocabulary
# sheet's M column:
#   - "DNBSEQ-T10x4RS" right before "DNBSEQ-T7"
#   - "Illumina NovaSeq X Plus" right before "Illumina iSeq 100"
# and grow the "instrumentmodel" named range accordingly (M1:M83 -> M1:M85).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cv_experiment")

# --- Step 1: insert "DNBSEQ-T10x4RS" before "DNBSEQ-T7" (currently row 30) ---
# Shift M30:M83 down to M31:M84 by copying values from the bottom up, then
# place the new value in the freed cell M30. We deliberately copy cell-by-cell
# (instead of using Range.Insert, which shifts the *entire* row and would
# also disturb the independent G/H/I columns on this sheet) so only column M
# is affected.
for ($r = 83; $r -ge 30; $r--) {
    $ws.Cells.Item($r + 1, 13).Value = $ws.Cells.Item($r, 13).Value2
}
$ws.Cells.Item(30, 13).Value = "DNBSEQ-T10x4RS"

# --- Step 2: insert "Illumina NovaSeq X Plus" before "Illumina iSeq 100" ---
# After step 1, "Illumina iSeq 100" moved from row 57 to row 58.
for ($r = 84; $r -ge 58; $r--) {
    $ws.Cells.Item($r + 1, 13).Value = $ws.Cells.Item($r, 13).Value2
}
$ws.Cells.Item(58, 13).Value = "Illumina NovaSeq X Plus"

# --- Step 3: grow the "instrumentmodel" defined name to cover the two new rows ---
$wb.Names.Item("instrumentmodel").RefersTo = "='cv_experiment'!`$M`$1:`$M`$85"
